$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.966.00'
$ws.Range("D3").Value = '2.294.30'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''299.97'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '''97.17'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '''0.505'
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").Value = '''33.69'
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").Value = '''0.0792'
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").Value = '''49.15'
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").Value = '''17.08'
$ws.Range("E14").Value = '  +11.82%  '
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '2.650.63'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '2.290.27'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").Value = '''0.806'
$ws.Range("E18").Value = '  +2.26%  '
$ws.Range("D19").Value = '42.911.60'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("D20").Value = '''11.67'
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").Value = '''67.48'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").Value = '''236.44'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("E25").Value = '  +5.33%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").Value = '''24.34'
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("D29").Value = '''2.17'
$ws.Range("E29").Value = '  -5.34%  '
$ws.Range("D30").Value = '''166.59'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").Value = '''33.77'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").Value = '''9.09'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("D36").Value = '''2.41'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '''16.73'
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("D38").Value = '''0.0695'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").Value = '''2.82'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").Value = '1.984.09'
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").Value = '''9.87'
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("D47").Value = '''17.43'
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").Value = '2.527.71'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("D50").Value = '''53.04'
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -2.18%  '
